$p = $ppt.ActivePresentation

# The speaker notes for slide 4 read "P = variabes (columnas)". In the
# source deck that sentence was typed as separate runs ("P = " / "variabes" /
# " " / "(columnas)"); re-typing the " (columnas)" tail merges the trailing
# " " run with the "(columnas)" run into a single " (columnas)" run (and
# leaves the cursor - hence the paragraph end mark - dirty). Locate that
# notes placeholder robustly (by content) and rewrite its text so the
# runs collapse the same way.

$targetText = "P = variabes (columnas)"
$found = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $notesPage = $slide.NotesPage
    for ($shi = 1; $shi -le $notesPage.Shapes.Count; $shi++) {
        $shape = $notesPage.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $targetText) {
                $tr.Text = $targetText
                $found = $true
            }
        }
    }
}

if (-not $found) {
    # Fallback to the known location (slide 4, notes body placeholder)
    $slide = $p.Slides.Item(4)
    $shape = $slide.NotesPage.Shapes.Item(2)
    $shape.TextFrame.TextRange.Text = $targetText
}
